# Insert a new "Industry" column before the existing "Mutual Fund" column (column C)
# and populate it with the industry classification for each holding row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing column C (Mutual Fund, Status, Jan_2026, Dec_2025, Oct_2025, MoM, QoQ)
# one place to the right, opening up a blank column C for "Industry".
$ws.Columns.Item(3).Insert()

# Header
$ws.Cells.Item(1, 3).Value = "Industry"

# Industry values per data row (row number -> industry name)
$industries = @{
    2  = "Banks"
    3  = "Aerospace & Defense"
    4  = "Banks"
    5  = "Petroleum Products"
    6  = "Aerospace & Defense"
    7  = "Power"
    8  = "Engineering Services"
    9  = "Aerospace & Defense"
    10 = "Banks"
    11 = "Power"
    12 = "Non - Ferrous Metals"
    13 = "Finance"
    14 = "Gas"
    15 = "Oil"
    16 = "Agricultural, Commercial & Construction Vehicles"
    17 = "Insurance"
    18 = "Industrial Manufacturing"
    19 = "Power"
    20 = "Petroleum Products"
    21 = "Power"
    22 = "Finance"
    23 = "Industrial Manufacturing"
    24 = "Leisure Services"
    25 = "Agricultural, Commercial & Construction Vehicles"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}
